$wb = $excel.ActiveWorkbook

# =========================================================================
# Add 8 new fish creatures to the "Creatures" sheet (rows 118-125) plus a
# matching "Sea Fish: Part 2" document entry on the "Documents" sheet.
# =========================================================================

$ws1  = $wb.Worksheets.Item("Creatures")
$ws11 = $wb.Worksheets.Item("Documents")

# ---- Creatures: stamp formatting for the new rows from the last row ----
$ws1.Range("A117:F117").Copy()
$ws1.Range("A118:F125").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$fishRows = @(
    @{Row=118; Name="Giant Anglerfish";   CR=7},
    @{Row=119; Name="Largetooth Sawfish"; CR=1},
    @{Row=120; Name="Smalltooth Sawfish"; CR=0.25},
    @{Row=121; Name="Giant Stargazer";    CR=6},
    @{Row=122; Name="Stonefish";          CR=0.5},
    @{Row=123; Name="Giant Stonefish";    CR=5},
    @{Row=124; Name="Swordfish";          CR=1},
    @{Row=125; Name="Great Swordfish";    CR=3}
)

# ---- fill in the plain text / number values first, so the new shared
#      strings land in the same order as the source names ----
foreach ($fr in $fishRows) {
    $r = $fr.Row
    $ws1.Range("A$r").Value = $fr.Name
    $ws1.Range("B$r").Value = $fr.CR
    $ws1.Range("C$r").Value = "Beast"
    $ws1.Range("E$r").Value = "Complete"
    $ws1.Range("F$r").Value = "Publicly Released"
}

# ---- Documents: add the "Sea Fish: Part 2" source document row (80) ----
$ws11.Range("A79:F79").Copy()
$ws11.Range("A80:F80").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws11.Range("B80").Value = "Complete"
$ws11.Range("C80").Value = "Yes"
$ws11.Range("D80").Value = "Publicly Released"
$ws11.Range("E80").Value = "Source Doc"

# ---- now add the hyperlinks (this also creates the "Sea Fish: Part 2"
#      shared string, last, matching the diff's string ordering) ----
$docUrl = "https://editor.gmbinder.com/documents/edit/-N9FseaFishPart2xxxx"

foreach ($fr in $fishRows) {
    $r = $fr.Row
    $ws1.Hyperlinks.Add($ws1.Range("D$r"), $docUrl, "", "", "Sea Fish: Part 2")
    # Hyperlinks.Add() stamps its own font style onto the cell; restore the
    # original column formatting that PasteSpecial already applied.
    $ws1.Range("D117").Copy()
    $ws1.Range("D$r").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

$ws11.Hyperlinks.Add($ws11.Range("A80"), $docUrl, "", "", "Sea Fish: Part 2")
$ws11.Range("A79").Copy()
$ws11.Range("A80").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# =========================================================================
# View-state bookkeeping to mirror the author's final selection/scroll
# position and active sheet.
# =========================================================================

$ws1.Range("D125").Select()
$ws11.Range("B81").Select()
$ws11.Activate()
